$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 230.33333
$ws.Range("I2").Value = 246.625
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 246.625
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = -133.625
$ws.Range("N2").Value = -326
$ws.Range("H38").Value = 3035.4
$ws.Range("I38").Value = 59
$ws.Range("J38").Value = 7500
$ws.Range("K38").Value = 177
$ws.Range("L38").Value = 22500
$ws.Range("M38").Value = 195
$ws.Range("N38").Value = -23244
$ws.Range("H39").Value = 390.75
$ws.Range("I39").Value = 325.2
$ws.Range("K39").Value = 975.5999999999999
$ws.Range("M39").Value = -679.5999999999999
$ws.Range("H40").Value = 3870.7
$ws.Range("I40").Value = 5533.6665
$ws.Range("K40").Value = 5533.6665
$ws.Range("M40").Value = -5358.6665
$ws.Range("H86").Value = 9336.727999999999
$ws.Range("I86").Value = 9088
$ws.Range("K86").Value = 9088
$ws.Range("M86").Value = -7965
$ws.Range("H89").Value = 9336.727999999999
$ws.Range("I89").Value = 9088
$ws.Range("K89").Value = 45440
$ws.Range("M89").Value = -39824
$ws.Range("H111").Value = 3051.75
$ws.Range("I111").Value = 5334.143
$ws.Range("K111").Value = 16002.429
$ws.Range("M111").Value = -12935.429
$ws.Range("H112").Value = 1523.5
$ws.Range("J112").Value = 1523.5
$ws.Range("L112").Value = 4570.5
$ws.Range("N112").Value = -6786.5
$ws.Range("I127").Value = 999.6667
$ws.Range("J127").Value = 999
$ws.Range("K127").Value = 2999.0001
$ws.Range("L127").Value = 2997
$ws.Range("M127").Value = 1960.9999
$ws.Range("N127").Value = -12917
$ws.Range("H138").Value = 3237.6667
$ws.Range("I138").Value = 754.1111
$ws.Range("K138").Value = 2262.3333
$ws.Range("M138").Value = 2877.6667
$ws.Range("H141").Value = 3547.5557
$ws.Range("I141").Value = 3821.8333
$ws.Range("J141").Value = 2999
$ws.Range("K141").Value = 11465.4999
$ws.Range("L141").Value = 8997
$ws.Range("M141").Value = -6285.499899999999
$ws.Range("N141").Value = -19357

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1517.3334
$ws.Range("I2").Value = 1677.5
$ws.Range("K2").Value = 1677.5
$ws.Range("M2").Value = -1564.5
$ws.Range("H25").Value = 2500
$ws.Range("I25").Value = 2500
$ws.Range("K25").Value = 2500
$ws.Range("M25").Value = -2098
$ws.Range("H32").Value = 3531714
$ws.Range("I32").Value = 1688231.5
$ws.Range("K32").Value = 1688231.5
$ws.Range("M32").Value = -1687944.5
$ws.Range("H35").Value = 400
$ws.Range("I35").Value = 400
$ws.Range("K35").Value = 400
$ws.Range("M35").Value = 6
$ws.Range("H45").Value = 28401.125
$ws.Range("I45").Value = 28401.125
$ws.Range("K45").Value = 28401.125
$ws.Range("M45").Value = -28024.125
$ws.Range("H61").Value = 1966.4482
$ws.Range("I61").Value = 1251.65
$ws.Range("K61").Value = 1251.65
$ws.Range("M61").Value = -1039.65
$ws.Range("H116").Value = 1517.3334
$ws.Range("I116").Value = 1677.5
$ws.Range("K116").Value = 1677.5
$ws.Range("M116").Value = 616.5
$ws.Range("H136").Value = 1966.4482
$ws.Range("I136").Value = 1251.65
$ws.Range("K136").Value = 3754.95
$ws.Range("M136").Value = -1204.95

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1517.3334
$ws.Range("I3").Value = 1677.5
$ws.Range("K3").Value = 1677.5
$ws.Range("M3").Value = -1563.5
$ws.Range("H134").Value = 1956.3334
$ws.Range("I134").Value = 1547.65
$ws.Range("K134").Value = 4642.950000000001
$ws.Range("M134").Value = -2107.950000000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1437
$ws.Range("I16").Value = 1405.5
$ws.Range("K16").Value = 1405.5
$ws.Range("M16").Value = -1118.5
$ws.Range("H31").Value = 5438294
$ws.Range("I31").Value = 2814.7334
$ws.Range("J31").Value = 15629817
$ws.Range("K31").Value = 2814.7334
$ws.Range("L31").Value = 15629817
$ws.Range("M31").Value = -2519.7334
$ws.Range("N31").Value = -15630407
$ws.Range("H34").Value = 5438294
$ws.Range("I34").Value = 2814.7334
$ws.Range("J34").Value = 15629817
$ws.Range("K34").Value = 2814.7334
$ws.Range("L34").Value = 15629817
$ws.Range("M34").Value = -2612.7334
$ws.Range("N34").Value = -15630221
$ws.Range("H58").Value = 1934.091
$ws.Range("I58").Value = 1376.5
$ws.Range("K58").Value = 1376.5
$ws.Range("M58").Value = -1173.5
$ws.Range("H113").Value = 1437
$ws.Range("I113").Value = 1405.5
$ws.Range("K113").Value = 1405.5
$ws.Range("M113").Value = 764.5
$ws.Range("H132").Value = 4855.643
$ws.Range("I132").Value = 4282.95
$ws.Range("J132").Value = 6287.375
$ws.Range("K132").Value = 12848.85
$ws.Range("L132").Value = 18862.125
$ws.Range("M132").Value = -10318.85
$ws.Range("N132").Value = -23922.125
$ws.Range("H134").Value = 5035.1665
$ws.Range("I134").Value = 5159.2856
$ws.Range("J134").Value = 4166.3335
$ws.Range("K134").Value = 15477.8568
$ws.Range("L134").Value = 12499.0005
$ws.Range("M134").Value = -12942.8568
$ws.Range("N134").Value = -17569.0005
$ws.Range("H136").Value = 1934.091
$ws.Range("I136").Value = 1376.5
$ws.Range("K136").Value = 4129.5
$ws.Range("M136").Value = -1579.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 108688.02
$ws.Range("I4").Value = 108344.92
$ws.Range("J4").Value = 125500
$ws.Range("K4").Value = 325034.76
$ws.Range("L4").Value = 376500
$ws.Range("M4").Value = -324922.76
$ws.Range("N4").Value = -376724
$ws.Range("H113").Value = 1473
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 1464
$ws.Range("K113").Value = 4500
$ws.Range("L113").Value = 4392
$ws.Range("M113").Value = -2330
$ws.Range("N113").Value = -8732
$ws.Range("H129").Value = 63139.223
$ws.Range("J129").Value = 63139.223
$ws.Range("L129").Value = 189417.669
$ws.Range("N129").Value = -199417.669
$ws.Range("H132").Value = 1604.4286
$ws.Range("I132").Value = 1280.3334
$ws.Range("K132").Value = 11523.0006
$ws.Range("M132").Value = -8993.000599999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 35000
$ws.Range("J49").Value = 35000
$ws.Range("L49").Value = 35000
$ws.Range("N49").Value = -35368
$ws.Range("H102").Value = 4886.552
$ws.Range("I102").Value = 1733.6666
$ws.Range("J102").Value = 5250.346
$ws.Range("K102").Value = 1733.6666
$ws.Range("L102").Value = 5250.346
$ws.Range("M102").Value = -111.6666
$ws.Range("N102").Value = -8494.346
$ws.Range("H122").Value = 3222.36
$ws.Range("I122").Value = 2398.1
$ws.Range("J122").Value = 3771.8667
$ws.Range("K122").Value = 7194.299999999999
$ws.Range("L122").Value = 11315.6001
$ws.Range("M122").Value = -4744.299999999999
$ws.Range("N122").Value = -16215.6001
$ws.Range("H123").Value = 74166.336
$ws.Range("I123").Value = 69999
$ws.Range("J123").Value = 76250
$ws.Range("K123").Value = 69999
$ws.Range("L123").Value = 76250
$ws.Range("M123").Value = -67549
$ws.Range("N123").Value = -81150
$ws.Range("H126").Value = 9088.450000000001
$ws.Range("I126").Value = 2073
$ws.Range("K126").Value = 6219
$ws.Range("M126").Value = -3749
$ws.Range("H132").Value = 2334.3635
$ws.Range("I132").Value = 2009.75
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 6029.25
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -3499.25
$ws.Range("N132").Value = -14660

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3000
$ws.Range("I22").Value = 3000
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -2705
$ws.Range("N22").Value = -3590
$ws.Range("H27").Value = 3000
$ws.Range("I27").Value = 3000
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -2893
$ws.Range("N27").Value = -3214
$ws.Range("H40").Value = 34334
$ws.Range("I40").Value = 40000.8
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 40000.8
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -39864.8
$ws.Range("N40").Value = -6272
$ws.Range("H127").Value = 66213
$ws.Range("J127").Value = 66213
$ws.Range("L127").Value = 66213
$ws.Range("N127").Value = -76133
$ws.Range("H131").Value = 54990
$ws.Range("J131").Value = 54990
$ws.Range("L131").Value = 54990
$ws.Range("N131").Value = -65070
$ws.Range("H132").Value = 5009.5386
$ws.Range("I132").Value = 5370.125
$ws.Range("J132").Value = 4432.6
$ws.Range("K132").Value = 16110.375
$ws.Range("L132").Value = 13297.8
$ws.Range("M132").Value = -13580.375
$ws.Range("N132").Value = -18357.8

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3890.9583
$ws.Range("I81").Value = 4582.75
$ws.Range("J81").Value = 3199.1667
$ws.Range("K81").Value = 9165.5
$ws.Range("L81").Value = 6398.3334
$ws.Range("M81").Value = -8104.5
$ws.Range("N81").Value = -8520.3334
$ws.Range("H84").Value = 3890.9583
$ws.Range("I84").Value = 4582.75
$ws.Range("J84").Value = 3199.1667
$ws.Range("K84").Value = 45827.5
$ws.Range("L84").Value = 31991.667
$ws.Range("M84").Value = -40523.5
$ws.Range("N84").Value = -42599.667
$ws.Range("H122").Value = 17859264
$ws.Range("I122").Value = 2265.8333
$ws.Range("K122").Value = 6797.499899999999
$ws.Range("M122").Value = -4347.499899999999
$ws.Range("H126").Value = 24000
$ws.Range("I126").Value = 24000
$ws.Range("K126").Value = 72000
$ws.Range("M126").Value = -69530
$ws.Range("H132").Value = 4440.067
$ws.Range("I132").Value = 4522.4287
$ws.Range("K132").Value = 13567.2861
$ws.Range("M132").Value = -11037.2861
